$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 30
$ws.Range("F3").Value = 491
$ws.Range("F5").Value = 91
$ws.Range("F6").Value = 335
$ws.Range("F7").Value = 1314
$ws.Range("F10").Value = 1351
$ws.Range("F13").Value = 189
$ws.Range("F17").Value = 1702
$ws.Range("F18").Value = 630
$ws.Range("F19").Value = 274
$ws.Range("F20").Value = 292
$ws.Range("F21").Value = 3216
$ws.Range("F22").Value = 28
$ws.Range("F23").Value = 413
$ws.Range("F24").Value = 940
$ws.Range("F25").Value = 1221
$ws.Range("F27").Value = 2862
$ws.Range("F28").Value = 1661
$ws.Range("F29").Value = 90
$ws.Range("F30").Value = 127
$ws.Range("F31").Value = 673
$ws.Range("F32").Value = 876
$ws.Range("F33").Value = 7
$ws.Range("F34").Value = 1939
$ws.Range("F35").Value = 911
$ws.Range("F36").Value = 1944
$ws.Range("F37").Value = 212
$ws.Range("F38").Value = 364
$ws.Range("F39").Value = 105
$ws.Range("F40").Value = 851
$ws.Range("F41").Value = 51
$ws.Range("F42").Value = 915
$ws.Range("F43").Value = 817
$ws.Range("F44").Value = 1060
$ws.Range("F45").Value = 147
$ws.Range("F48").Value = 238
$ws.Range("F49").Value = 3384

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 816
$ws.Range("F14").Value = 26

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 491
$ws.Range("F4").Value = 91
$ws.Range("F7").Value = 335
$ws.Range("F8").Value = 1314
$ws.Range("F11").Value = 1351
$ws.Range("F14").Value = 189
$ws.Range("F17").Value = 1702
$ws.Range("F18").Value = 630
$ws.Range("F19").Value = 274
$ws.Range("F20").Value = 292
$ws.Range("F21").Value = 3216
$ws.Range("F22").Value = 28
$ws.Range("F23").Value = 413
$ws.Range("F25").Value = 1221
$ws.Range("F26").Value = 2862
$ws.Range("F27").Value = 1661
$ws.Range("F28").Value = 90
$ws.Range("F29").Value = 127
$ws.Range("F30").Value = 816
$ws.Range("F31").Value = 26
$ws.Range("F32").Value = 876
$ws.Range("F33").Value = 1939
$ws.Range("F35").Value = 911
$ws.Range("F36").Value = 1944
$ws.Range("F37").Value = 364
$ws.Range("F38").Value = 105
$ws.Range("F39").Value = 851
$ws.Range("F40").Value = 915
$ws.Range("F41").Value = 817
$ws.Range("F42").Value = 1060
$ws.Range("F43").Value = 147
$ws.Range("F47").Value = 238
$ws.Range("F48").Value = 3384
